$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, derived from the commit diff.
$updates = [ordered]@{
    'D2' = '302.46'
    'E2' = '-1.22%'
    'G2' = '23'
    'D3' = '35.29'
    'E3' = '-2.70%'
    'G3' = '23'
    'D4' = '5.042'
    'E4' = '-0.54%'
    'G4' = '23'
    'D5' = '0.07894'
    'E5' = '-0.63%'
    'G5' = '23'
    'D6' = '1.931'
    'E6' = '-11.14%'
    'G6' = '23'
    'B7' = 'KuCoinToken'
    'C7' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D7' = '7.831'
    'E7' = '-2.38%'
    'G7' = '23'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D8' = '2.846'
    'E8' = '8.14%'
    'G8' = '23'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '0.9280'
    'E9' = '-0.44%'
    'G9' = '23'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.1078'
    'E10' = '9.02%'
    'G10' = '23'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1890'
    'E11' = '0.88%'
    'G11' = '23'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.09317'
    'E12' = '3.24%'
    'G12' = '23'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '0.03681'
    'E13' = '1.49%'
    'G13' = '23'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '0.09935'
    'E14' = '0.01%'
    'G14' = '23'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '0.001446'
    'E15' = '0.95%'
    'G15' = '23'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.005739'
    'E16' = '1.48%'
    'G16' = '23'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.469'
    'E17' = '0.87%'
    'G17' = '23'
    'B18' = 'GateToken'
    'C18' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D18' = '4.115'
    'E18' = '-1.20%'
    'G18' = '23'
    'D19' = '0.3434'
    'E19' = '1.86%'
    'G19' = '23'
    'E20' = '-2.97%'
    'G20' = '23'
    'D21' = '5.117'
    'E21' = '0.64%'
    'G21' = '23'
    'E22' = '0.46%'
    'G22' = '23'
    'D23' = '0.04533'
    'E23' = '-1.16%'
    'G23' = '23'
    'D24' = '0.001230'
    'E24' = '-0.89%'
    'G24' = '23'
    'D25' = '0.004678'
    'E25' = '-1.84%'
    'G25' = '23'
    'D26' = '0.0001256'
    'E26' = '-3.61%'
    'G26' = '23'
    'D27' = '0.0004473'
    'E27' = '-5.60%'
    'G27' = '23'
    'G28' = '23'
    'G29' = '23'
    'G30' = '23'
    'G31' = '23'
    'G32' = '23'
    'G33' = '23'
    'G34' = '23'
    'G35' = '23'
    'G36' = '23'
    'G37' = '23'
    'G38' = '23'
    'D39' = '0.01887'
    'E39' = '-3.40%'
    'G39' = '23'
    'D40' = '0.04722'
    'E40' = '-4.12%'
    'G40' = '23'
    'D41' = '0.007625'
    'E41' = '-2.12%'
    'G41' = '23'
    'D42' = '0.01004'
    'E42' = '32.06%'
    'G42' = '23'
    'D43' = '0.1343'
    'E43' = '-3.64%'
    'G43' = '23'
    'D44' = '0.002125'
    'E44' = '1.13%'
    'G44' = '23'
    'D45' = '0.01127'
    'E45' = '0.16%'
    'G45' = '23'
    'D46' = '0.00006337'
    'E46' = '1.97%'
    'G46' = '23'
    'D47' = '0.00000000754'
    'E47' = '0.42%'
    'G47' = '23'
    'E48' = '23.41%'
    'G48' = '23'
    'D49' = '0.001307'
    'E49' = '-27.48%'
    'G49' = '23'
    'D50' = '0.00002112'
    'E50' = '0.42%'
    'G50' = '23'
    'D51' = '0.0002012'
    'E51' = '0.42%'
    'G51' = '23'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (prices, percents,
    # hour counters) are not reinterpreted as numbers by Excel.
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$addr]
    $cell.NumberFormat = 'General'
    $cell.ClearFormats()
}
